$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.477.88'
$ws.Range('E2').Value = '  +2.53%  '
$ws.Range('D3').Value = '2.351.72'
$ws.Range('E3').Value = '  +5.88%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.643'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.99%  '
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.637'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.89'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('E13').Value = '  +4.44%  '
$ws.Range('E14').Value = '  +2.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.36'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.58%  '
$ws.Range('D16').Value = '2.706.27'
$ws.Range('E16').Value = '  +6.03%  '
$ws.Range('D17').Value = '2.415.79'
$ws.Range('E17').Value = '  +8.94%  '
$ws.Range('D18').Value = '43.445.75'
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('E19').Value = '  +3.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '75.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.17%  '
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '257.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.29'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.29%  '
$ws.Range('B31').Value = 'WEMIXToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '173.46'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0931'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.49%  '
$ws.Range('E34').Value = '  +7.82%  '
$ws.Range('E35').Value = '  +5.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.98'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.17'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0377'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.105'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.47%  '
$ws.Range('E40').Value = '  +13.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('E42').Value = '  +13.90%  '
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.43%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +11.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.02'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.101'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.465'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.56%  '
